$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------
$RPR = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

function Get-RunXml($text, $preserve) {
    $space = ""
    if ($preserve) { $space = ' xml:space="preserve"' }
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    return "<w:r>$RPR<w:t$space>$escaped</w:t></w:r>"
}

function Get-PkgXml($bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        "<w:body>$bodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
}

# Builds a full new list-item paragraph (pStyle a5, ilvl=2, numId=4) whose
# text is the concatenation of one or more runs.
function Get-ListItemParaXml($runsXml) {
    $ppr = "<w:pPr><w:pStyle w:val=`"a5`"/><w:numPr><w:ilvl w:val=`"2`"/><w:numId w:val=`"4`"/></w:numPr><w:spacing w:line=`"360`" w:lineRule=`"auto`"/>$RPR</w:pPr>"
    return "<w:p>$ppr$runsXml</w:p>"
}

# Inserts a brand-new paragraph (with its own pPr+runs) right after the
# paragraph currently at 1-based index $afterIndex, returning the index of
# the freshly-created paragraph.
function Insert-NewParagraphAfter($afterIndex, $paraXml) {
    $prev = $d.Paragraphs.Item($afterIndex)
    $prev.Range.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $full = $d.Paragraphs.Item($newIndex).Range
    $full.InsertXML((Get-PkgXml $paraXml))
    return $newIndex
}

# ---------------------------------------------------------------------------
# 1) Merge the two runs of "If any of the words in the s" / "olution list
#    are 5 letter words" into a single run, dropping the _GoBack bookmark
#    that used to sit between them.
# ---------------------------------------------------------------------------
$p26 = $d.Paragraphs.Item(26)
$full26 = $p26.Range
$r26 = $d.Range($full26.Start, $full26.End - 1)
$runXml = Get-RunXml "If any of the words in the solution list are 5 letter words" $false
$r26.InsertXML((Get-PkgXml "<w:p>$runXml</w:p>"))

# ---------------------------------------------------------------------------
# 2) Re-add the _GoBack bookmark inside the (empty) paragraph that follows
#    "If any of the words in the solution list are repeated".
# ---------------------------------------------------------------------------
$p28 = $d.Paragraphs.Item(28)
$full28 = $p28.Range
$r28 = $d.Range($full28.Start, $full28.Start)
$bookmarkXml = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$r28.InsertXML((Get-PkgXml "<w:p>$bookmarkXml</w:p>"))

# ---------------------------------------------------------------------------
# 3) Fill in the (empty) "White box testing" sub-bullet with "Test all the
#    path", then append six new ilvl=2 bullets after it.
# ---------------------------------------------------------------------------
$p33 = $d.Paragraphs.Item(33)
$full33 = $p33.Range
$r33 = $d.Range($full33.Start, $full33.End - 1)
$testAllXml = Get-RunXml "Test all the path" $false
$r33.InsertXML((Get-PkgXml "<w:p>$testAllXml</w:p>"))

$idx = 33

# "Try one example which the starting word is not 5 word letter long"
$runs = Get-RunXml "Try one example which the starting word is not 5 word letter long" $false
$paraXml = Get-ListItemParaXml $runs
$idx = Insert-NewParagraphAfter $idx $paraXml

# "Try" / " one example " / "which" / " the " / "ending word is " / "not 5 word letter long"
$runs = (Get-RunXml "Try" $false) +
        (Get-RunXml " one example " $true) +
        (Get-RunXml "which" $false) +
        (Get-RunXml " the " $true) +
        (Get-RunXml "ending word is " $true) +
        (Get-RunXml "not 5 word letter long" $false)
$paraXml = Get-ListItemParaXml $runs
$idx = Insert-NewParagraphAfter $idx $paraXml

# "Try one example which the starting word is not in the dictionary"
$runs = Get-RunXml "Try one example which the starting word is not in the dictionary" $false
$paraXml = Get-ListItemParaXml $runs
$idx = Insert-NewParagraphAfter $idx $paraXml

# "Ending w" / "ord" / " " / "is" / " not in the dictionary"
$runs = (Get-RunXml "Ending w" $false) +
        (Get-RunXml "ord" $false) +
        (Get-RunXml " " $true) +
        (Get-RunXml "is" $false) +
        (Get-RunXml " not in the dictionary" $true)
$paraXml = Get-ListItemParaXml $runs
$idx = Insert-NewParagraphAfter $idx $paraXml

# "Try one example which doesn't have any solution" (curly apostrophe)
$curlyApos = [char]0x2019
$doesntText = "Try one example which doesn" + $curlyApos + "t have any solution"
$runs = Get-RunXml $doesntText $false
$paraXml = Get-ListItemParaXml $runs
$idx = Insert-NewParagraphAfter $idx $paraXml

# "Try one example which the ending word and starting words are same"
$runs = Get-RunXml "Try one example which the ending word and starting words are same" $false
$paraXml = Get-ListItemParaXml $runs
$idx = Insert-NewParagraphAfter $idx $paraXml

Write-Output "done; paragraph count = $($d.Paragraphs.Count)"
